$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H11").Value = "   "
$ws.Range("H10").Value = "  "
$ws.Range("H9").Value = "     "
$ws.Range("H14").Value = "   "
$ws.Range("H18").Value = "   "
$ws.Range("H20").Value = " "
$ws.Range("H23").Value = "  "
$ws.Range("H18").Select() | Out-Null
